$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.726.00"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "3.437.86"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'582.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'172.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.601"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").Value = "3.435.24"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'0.130"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("E12").Value = "  -2.76%  "
$ws.Range("D13").Value = "4.032.92"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("D15").Value = "'28.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.41%  "
$ws.Range("D16").Value = "65.795.93"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "3.448.44"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").Value = "'13.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "'366.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("D23").Value = "'71.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").Value = "'0.998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").Value = "'0.530"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("E26").Value = "  +3.21%  "
$ws.Range("D27").Value = "'9.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("E28").Value = "  +2.70%  "
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "'23.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'5.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  -4.73%  "
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").Value = "'160.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "'0.877"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").Value = "'28.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.84%  "
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("D42").Value = "2.751.72"
$ws.Range("E42").Value = "  +2.89%  "
$ws.Range("D43").Value = "'6.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("D46").Value = "'40.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").Value = "'24.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").Value = "'325.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "'6.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.68%  "
